$d = $word.ActiveDocument

# Locate the heading paragraph "Lista de ítem con la nomenclatura" -- the new
# content must be inserted immediately before it (right after the preceding
# heading "Definición de la nomenclatura de ítem").
$findRng = $d.Range(0, 0)
$found = $findRng.Find.Execute("Lista de " + [char]237 + "tem con la nomenclatura", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target heading paragraph"
}

$insStart = $findRng.Start

# Build the three new paragraphs as a raw OOXML fragment and insert them via
# InsertXML so that paragraph/run formatting (indentation, justification,
# language, bold, proofErr markers, list style) is reproduced exactly as in
# the target revision. The fragment ends with an "open" empty paragraph,
# which Word merges into the following (pre-existing) heading paragraph
# rather than creating an extra empty paragraph.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/part.xml" pkg:contentType="application/xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:left="1440"/><w:jc w:val="both"/><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>En este caso, la nomenclatura se ha definido de acuerdo con el acr&#243;nimo del nombre del proyecto (</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>PARKING</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> SOFT) junto con el acr&#243;nimo del elemento o &#237;tem que se manejar&#225;. Siendo as&#237;, la nomenclatura de &#237;tem tendr&#225; la siguiente estructura:</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="3024" w:firstLine="576"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:lang w:val="es-MX"/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="es-MX"/></w:rPr><w:t>PS - &#8220;Acr&#243;nimo del &#237;tem&#8221;</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:ind w:left="504"/></w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRng = $d.Range($insStart, $insStart)
$insertRng.InsertXML($xml)
